$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 800
$ws.Range("I18").Value = 800
$ws.Range("K18").Value = 800
$ws.Range("M18").Value = -516
$ws.Range("H41").Value = 1942.6666
$ws.Range("I41").Value = 1715.5454
$ws.Range("K41").Value = 1715.5454
$ws.Range("M41").Value = -1275.5454
$ws.Range("H55").Value = 384.63635
$ws.Range("I55").Value = 206.4
$ws.Range("K55").Value = 206.4
$ws.Range("M55").Value = 7.599999999999994
$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 5400
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 5400
$ws.Range("M74").Value = -2064
$ws.Range("N74").Value = -7272
$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 5400
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 27000
$ws.Range("M77").Value = -10320
$ws.Range("N77").Value = -36360
$ws.Range("H98").Value = 904.7646999999999
$ws.Range("I98").Value = 904.7646999999999
$ws.Range("K98").Value = 904.7646999999999
$ws.Range("M98").Value = 593.2353000000001
$ws.Range("H106").Value = 9743.923000000001
$ws.Range("I106").Value = 1853
$ws.Range("K106").Value = 1853
$ws.Range("M106").Value = -1222
$ws.Range("H111").Value = 475
$ws.Range("J111").Value = 550
$ws.Range("L111").Value = 1650
$ws.Range("N111").Value = -7784
$ws.Range("H122").Value = 904.7646999999999
$ws.Range("I122").Value = 904.7646999999999
$ws.Range("K122").Value = 2714.2941
$ws.Range("M122").Value = -264.2941000000001
$ws.Range("H132").Value = 1313.2391
$ws.Range("I132").Value = 893.0952
$ws.Range("K132").Value = 2679.2856
$ws.Range("M132").Value = -149.2856000000002
$ws.Range("H137").Value = 2387.8948
$ws.Range("I137").Value = 2038.9656
$ws.Range("J137").Value = 3512.2222
$ws.Range("K137").Value = 6116.8968
$ws.Range("L137").Value = 10536.6666
$ws.Range("M137").Value = -3566.8968
$ws.Range("N137").Value = -15636.6666
$ws.Range("H138").Value = 2771
$ws.Range("I138").Value = 1949.3636
$ws.Range("J138").Value = 3181.818
$ws.Range("K138").Value = 5848.0908
$ws.Range("L138").Value = 9545.454000000002
$ws.Range("M138").Value = -708.0907999999999
$ws.Range("N138").Value = -19825.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 407060.9
$ws.Range("I45").Value = 633949.75
$ws.Range("K45").Value = 633949.75
$ws.Range("M45").Value = -633572.75
$ws.Range("H132").Value = 1958.9318
$ws.Range("I132").Value = 1073.8379
$ws.Range("K132").Value = 3221.5137
$ws.Range("M132").Value = -691.5137
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9743.931
$ws.Range("I20").Value = 8345.235000000001
$ws.Range("J20").Value = 11725.417
$ws.Range("K20").Value = 8345.235000000001
$ws.Range("L20").Value = 11725.417
$ws.Range("M20").Value = -8098.235000000001
$ws.Range("N20").Value = -12219.417
$ws.Range("H86").Value = 2468.2222
$ws.Range("I86").Value = 1600.1333
$ws.Range("J86").Value = 3088.2856
$ws.Range("K86").Value = 1600.1333
$ws.Range("L86").Value = 3088.2856
$ws.Range("M86").Value = -477.1333
$ws.Range("N86").Value = -5334.2856
$ws.Range("H89").Value = 2468.2222
$ws.Range("I89").Value = 1600.1333
$ws.Range("J89").Value = 3088.2856
$ws.Range("K89").Value = 8000.666499999999
$ws.Range("L89").Value = 15441.428
$ws.Range("M89").Value = -2384.666499999999
$ws.Range("N89").Value = -26673.428
$ws.Range("H129").Value = 66662.336
$ws.Range("J129").Value = 66662.336
$ws.Range("L129").Value = 66662.336
$ws.Range("N129").Value = -76662.336
$ws.Range("H134").Value = 1327.7428
$ws.Range("I134").Value = 1318
$ws.Range("J134").Value = 1431.6666
$ws.Range("K134").Value = 3954
$ws.Range("L134").Value = 4294.9998
$ws.Range("M134").Value = -1419
$ws.Range("N134").Value = -9364.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1332.2858
$ws.Range("I16").Value = 1025.9333
$ws.Range("J16").Value = 2098.1667
$ws.Range("K16").Value = 1025.9333
$ws.Range("L16").Value = 2098.1667
$ws.Range("M16").Value = -738.9332999999999
$ws.Range("N16").Value = -2672.1667
$ws.Range("H31").Value = 3538.0833
$ws.Range("I31").Value = 1617.25
$ws.Range("J31").Value = 7379.75
$ws.Range("K31").Value = 1617.25
$ws.Range("L31").Value = 7379.75
$ws.Range("M31").Value = -1322.25
$ws.Range("N31").Value = -7969.75
$ws.Range("H34").Value = 3538.0833
$ws.Range("I34").Value = 1617.25
$ws.Range("J34").Value = 7379.75
$ws.Range("K34").Value = 1617.25
$ws.Range("L34").Value = 7379.75
$ws.Range("M34").Value = -1415.25
$ws.Range("N34").Value = -7783.75
$ws.Range("H105").Value = 1541.091
$ws.Range("I105").Value = 1495.2
$ws.Range("K105").Value = 1495.2
$ws.Range("M105").Value = 251.8
$ws.Range("H107").Value = 1097.84
$ws.Range("I107").Value = 772.4
$ws.Range("J107").Value = 1586
$ws.Range("K107").Value = 772.4
$ws.Range("L107").Value = 1586
$ws.Range("M107").Value = 1147.6
$ws.Range("N107").Value = -5426
$ws.Range("H113").Value = 1332.2858
$ws.Range("I113").Value = 1025.9333
$ws.Range("J113").Value = 2098.1667
$ws.Range("K113").Value = 1025.9333
$ws.Range("L113").Value = 2098.1667
$ws.Range("M113").Value = 1144.0667
$ws.Range("N113").Value = -6438.1667
$ws.Range("H132").Value = 1895.9032
$ws.Range("I132").Value = 1819.8276
$ws.Range("K132").Value = 5459.4828
$ws.Range("M132").Value = -2929.4828

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4743.625
$ws.Range("J39").Value = 5574.8335
$ws.Range("L39").Value = 16724.5005
$ws.Range("N39").Value = -17312.5005
$ws.Range("H61").Value = 204.47058
$ws.Range("I61").Value = 120.44444
$ws.Range("K61").Value = 361.33332
$ws.Range("M61").Value = -146.33332
$ws.Range("H122").Value = 1319.1428
$ws.Range("J122").Value = 653.75
$ws.Range("L122").Value = 5883.75
$ws.Range("N122").Value = -10783.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5829.375
$ws.Range("I70").Value = 5806
$ws.Range("K70").Value = 5806
$ws.Range("M70").Value = -5536
$ws.Range("H73").Value = 5829.375
$ws.Range("I73").Value = 5806
$ws.Range("K73").Value = 5806
$ws.Range("M73").Value = -4870
$ws.Range("H80").Value = 4161.636
$ws.Range("J80").Value = 5796
$ws.Range("L80").Value = 5796
$ws.Range("N80").Value = -7792
$ws.Range("H83").Value = 4161.636
$ws.Range("J83").Value = 5796
$ws.Range("L83").Value = 28980
$ws.Range("N83").Value = -38964
$ws.Range("H113").Value = 2584.818
$ws.Range("I113").Value = 2789.6667
$ws.Range("K113").Value = 2789.6667
$ws.Range("M113").Value = -619.6667000000002
$ws.Range("H132").Value = 1979.2258
$ws.Range("I132").Value = 1830.28
$ws.Range("K132").Value = 5490.84
$ws.Range("M132").Value = -2960.84

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1135.5
$ws.Range("I82").Value = 1181.5
$ws.Range("J82").Value = 974.5
$ws.Range("K82").Value = 1181.5
$ws.Range("L82").Value = 974.5
$ws.Range("M82").Value = -820.5
$ws.Range("N82").Value = -1696.5
$ws.Range("H85").Value = 1135.5
$ws.Range("I85").Value = 1181.5
$ws.Range("J85").Value = 974.5
$ws.Range("K85").Value = 1181.5
$ws.Range("L85").Value = 974.5
$ws.Range("M85").Value = 66.5
$ws.Range("N85").Value = -3470.5
$ws.Range("H100").Value = 7727.278
$ws.Range("I100").Value = 7952.9287
$ws.Range("K100").Value = 7952.9287
$ws.Range("M100").Value = -7411.9287
$ws.Range("H132").Value = 4229.64
$ws.Range("I132").Value = 3410.4285
$ws.Range("J132").Value = 5272.273
$ws.Range("K132").Value = 10231.2855
$ws.Range("L132").Value = 15816.819
$ws.Range("M132").Value = -7701.2855
$ws.Range("N132").Value = -20876.819
$ws.Range("H136").Value = 2270.0625
$ws.Range("I136").Value = 1888.4
$ws.Range("K136").Value = 5665.200000000001
$ws.Range("M136").Value = -3115.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3711.9
$ws.Range("I100").Value = 4035.4443
$ws.Range("K100").Value = 8070.8886
$ws.Range("M100").Value = -7529.8886
$ws.Range("H122").Value = 2085.4285
$ws.Range("I122").Value = 3351.5
$ws.Range("J122").Value = 1579
$ws.Range("K122").Value = 10054.5
$ws.Range("L122").Value = 4737
$ws.Range("M122").Value = -7604.5
$ws.Range("N122").Value = -9637
$ws.Range("H126").Value = 4552.1577
$ws.Range("I126").Value = 3249.3572
$ws.Range("J126").Value = 8200
$ws.Range("K126").Value = 9748.071599999999
$ws.Range("L126").Value = 24600
$ws.Range("M126").Value = -7278.071599999999
$ws.Range("N126").Value = -29540
$ws.Range("H132").Value = 8296.277
$ws.Range("I132").Value = 8296.277
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 24888.831
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -22358.831
$ws.Range("H133").Value = 65355
$ws.Range("J133").Value = 65355
$ws.Range("L133").Value = 65355
$ws.Range("N133").Value = -75475
$ws.Range("H135").Value = 76196.8
$ws.Range("J135").Value = 76196.8
$ws.Range("L135").Value = 76196.8
$ws.Range("N135").Value = -86336.8
$ws.Range("H141").Value = 94000
$ws.Range("J141").Value = 94000
$ws.Range("L141").Value = 94000
$ws.Range("N141").Value = -104360
